$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the 9th submission row (50-xgboost, non-overlapping 3in1, preprocessed, 2 valid sets)
$ws.Range("A10").Value = "9_291115_0944_50_xgboost_with_non_overlap_3in1_preprocess_valid1_valid2_"
$ws.Range("B10").Value = 0.616
$ws.Range("C10").Value = "ensembled 50 xgboost, in non_overlapping 3in1 data set with features preprocessed, with 2 valid sets"

# Move the active selection down to the next empty row, as happens after typing a row of data
$ws.Range("C11").Select() | Out-Null
